$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

$ws.Range("D2").Value = "29.428.98"
$ws.Range("E2").Value = "  +0.01%  "
$ws.Range("D3").Value = "1.848.61"
$ws.Range("E3").Value = "  -0.08%  "
$ws.Range("D4").Value = "0.9986"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "240.74"
$ws.Range("E5").Value = "  -0.98%  "
$ws.Range("D6").Value = "0.6324"
$ws.Range("D7").Value = "0.9998"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "3.221.10"
$ws.Range("E8").Value = "  +74.37%  "
$ws.Range("D9").Value = "0.07589"
$ws.Range("E9").Value = "  +1.52%  "
$ws.Range("D10").Value = "0.2971"
$ws.Range("E10").Value = "  -0.75%  "
$ws.Range("D11").Value = "24.61"
$ws.Range("E11").Value = "  +1.15%  "
$ws.Range("E12").Value = "  +1.11%  "
$ws.Range("D13").Value = "4.990"
$ws.Range("E13").Value = "  -0.49%  "
$ws.Range("D14").Value = "0.6854"
$ws.Range("E14").Value = "  +0.20%  "
$ws.Range("D15").Value = "0.000009998"
$ws.Range("E15").Value = "  +5.38%  "
$ws.Range("D16").Value = "82.87"
$ws.Range("E16").Value = "  -0.93%  "
$ws.Range("D17").Value = "6.192"
$ws.Range("E17").Value = "  +0.95%  "
$ws.Range("D18").Value = "29.454.89"
$ws.Range("E18").Value = "  +0.05%  "
$ws.Range("D19").Value = "232.26"
$ws.Range("E19").Value = "  -2.27%  "
$ws.Range("E20").Value = "  -0.38%  "
$ws.Range("D21").Value = "0.9998"
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("E22").Value = "  -0.84%  "
$ws.Range("D23").Value = "0.9994"
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("D24").Value = "155.09"
$ws.Range("E24").Value = "  -1.08%  "
$ws.Range("D25").Value = "0.1390"
$ws.Range("E25").Value = "  -2.41%  "
$ws.Range("D26").Value = "8.442"
$ws.Range("E26").Value = "  -0.43%  "
$ws.Range("D27").Value = "17.67"
$ws.Range("E27").Value = "  -0.69%  "
$ws.Range("D28").Value = "1.473"
$ws.Range("E28").Value = "  -1.00%  "
$ws.Range("D29").Value = "0.05805"
$ws.Range("E29").Value = "  -3.93%  "
$ws.Range("D30").Value = "1.258"
$ws.Range("E30").Value = "  +0.77%  "
$ws.Range("D31").Value = "4.126"
$ws.Range("E31").Value = "  -0.23%  "
$ws.Range("D32").Value = "4.025"
$ws.Range("E32").Value = "  -1.17%  "
$ws.Range("D33").Value = "3.317.07"
$ws.Range("E33").Value = "  +65.49%  "
$ws.Range("D34").Value = "1.871"
$ws.Range("E34").Value = "  +0.72%  "
$ws.Range("D35").Value = "1.159"
$ws.Range("E35").Value = "  -1.54%  "
$ws.Range("D36").Value = "0.7224"
$ws.Range("E36").Value = "  +0.09%  "
$ws.Range("D37").Value = "2.594"
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("D38").Value = "1.249.82"
$ws.Range("E38").Value = "  +4.30%  "
$ws.Range("D39").Value = "2.794"
$ws.Range("E39").Value = "  +0.13%  "
$ws.Range("D40").Value = "0.01807"
$ws.Range("E40").Value = "  +1.42%  "
$ws.Range("D41").Value = "0.9025"
$ws.Range("E41").Value = "  -1.07%  "
$ws.Range("E42").Value = "  -2.02%  "
$ws.Range("D43").Value = "0.9990"
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("D44").Value = "101.45"
$ws.Range("E44").Value = "  -0.29%  "
$ws.Range("E45").Value = "  +0.92%  "
$ws.Range("D46").Value = "7.330"
$ws.Range("E46").Value = "  -1.64%  "
$ws.Range("D47").Value = "9.197"
$ws.Range("E47").Value = "  +1.49%  "
$ws.Range("D48").Value = "0.4018"
$ws.Range("E48").Value = "  -0.73%  "
$ws.Range("D49").Value = "1.695"
$ws.Range("E49").Value = "  +2.36%  "
$ws.Range("E50").Value = "  +0.06%  "
$ws.Range("E51").Value = "  +0.29%  "
